$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")
$ws.Activate() | Out-Null

# Cell style updates: J5 and L6 flip from the red "not done" fill to the
# green "done" fill already used by I5 / K6 in the same rows.
$ws.Range("J5").Interior.Color = $ws.Range("I5").Interior.Color
$ws.Range("L6").Interior.Color = $ws.Range("K6").Interior.Color

# Two more boxes ticked in the score grid.
$ws.Range("J11").Value = 1
$ws.Range("L12").Value = 1

# View changes: scroll the window right (B1 -> G1 as top-left) and move the
# selection to L6.
$ws.Range("G1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L6").Select() | Out-Null
